# Weekly update: two new "Arveja Verde" price records were reported for
# Comercializadora del Agro de Limarí and inserted into the historical
# series (sorted roughly by report date), pushing the existing rows down.
#
# Net effect vs. the original sheet:
#   - a new row is inserted at row 5 (date 2021-09-23)
#   - a new row is inserted at what becomes row 13 (date 2021-09-30)
#   - every other data row keeps its original values but shifts down by
#     one (rows 5-11 -> 6-12) or two (rows 12-34 -> 14-36; row 35 -> 37)
#   - the sheet dimension grows from A1:R35 to A1:R37

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 5 -----------------------------------
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value  = 2
$ws.Cells.Item(5, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(5, 3).Value  = "Coquimbo"
$ws.Cells.Item(5, 4).Value  = 44462
$ws.Cells.Item(5, 5).Value  = 4
$ws.Cells.Item(5, 6).Value  = 100112022
$ws.Cells.Item(5, 7).Value  = "Arveja Verde"
$ws.Cells.Item(5, 8).Value  = "Perfection"
$ws.Cells.Item(5, 9).Value  = "Primera"
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 22000
$ws.Cells.Item(5, 12).Value = 23000
$ws.Cells.Item(5, 13).Value = 22500
$ws.Cells.Item(5, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 900
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# --- Insert second new row at position 13 (after the shift above) --------
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value  = 2
$ws.Cells.Item(13, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13, 3).Value  = "Coquimbo"
$ws.Cells.Item(13, 4).Value  = 44469
$ws.Cells.Item(13, 5).Value  = 4
$ws.Cells.Item(13, 6).Value  = 100112022
$ws.Cells.Item(13, 7).Value  = "Arveja Verde"
$ws.Cells.Item(13, 8).Value  = "Perfection"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 600
$ws.Cells.Item(13, 11).Value = 22000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 23000
$ws.Cells.Item(13, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 920
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
